$d = $word.ActiveDocument

function Info($tag) {
  Write-Output "=== $tag ==="
  Write-Output "paras=$($d.Paragraphs.Count)"
  for ($i=1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    Write-Output "P$i [$($p.Range.Start)-$($p.Range.End)] text=[$($p.Range.Text)]"
  }
}

# --- Step 1: paragraph 1 run split: "git " / "dfsfdsfas" ---
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()
$p2 = $d.Paragraphs.Item(2)
$p2.Range.InsertAfter("dfsfdsfas")
$p1 = $d.Paragraphs.Item(1)
$markpos = $p1.Range.End - 1
$d.Range($markpos, $markpos+1).Delete()

# --- Step 2: two empty paragraphs ---
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()

$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()

# --- Step 3: "Fsdfsdaf" paragraph ---
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()
$lastp = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastp.Range.InsertAfter("Fsdfsdaf")

# --- Step 4: empty paragraph ---
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()

# --- Step 5: new paragraph for "fasd" -- insert first char "f" normally to avoid empty-para edge bug ---
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()
$lastp = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastp.Range.InsertAfter("f")
Info "structure with just f in last para"

$lastp = $d.Paragraphs.Item($d.Paragraphs.Count)
Write-Output "lastp range: $($lastp.Range.Start)-$($lastp.Range.End)"
$safePos = $lastp.Range.Start
Write-Output "safePos=$safePos"

$d.Bookmarks.Item("_GoBack").Delete()
$bm = $d.Bookmarks.Add("_GoBack", $d.Range($safePos, $safePos))
Write-Output "new bm: $($d.Bookmarks.Item('_GoBack').Start)-$($d.Bookmarks.Item('_GoBack').End)"

# nudge: insert "asd" via bookmark range so it rides to true end (note: bookmark is BEFORE the "f", so inserting after will push past "f" too? let's check)
$r = $d.Bookmarks.Item("_GoBack").Range
$r.InsertAfter("asd")
Info "final structure"
Write-Output "final bm: $($d.Bookmarks.Item('_GoBack').Start)-$($d.Bookmarks.Item('_GoBack').End)"
